$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsRuntimes = $wb.Worksheets.Item("Runtimes")

# --- Data sheet: new rows 710-720 ---
# row 710
$wsData.Cells.Item(710, 1).Value = 100
$wsData.Cells.Item(710, 2).Value = 5
$wsData.Cells.Item(710, 3).Value = 0.5
$wsData.Cells.Item(710, 4).Value = "ba-no-cycle"
$wsData.Cells.Item(710, 5).Value = "2021-07-07 21:57:08.700503"
$wsData.Cells.Item(710, 6).Value = 1
$wsData.Cells.Item(710, 7).Value = 0
$wsData.Cells.Item(710, 8).Value = 0
$wsData.Cells.Item(710, 9).Value = 0
$wsData.Cells.Item(710, 10).Value = 5
$wsData.Cells.Item(710, 11).Value = 5
$wsData.Cells.Item(710, 12).Value = 5
$wsData.Cells.Item(710, 13).Value = "-"
$wsData.Cells.Item(710, 14).Value = "-"
$wsData.Cells.Item(710, 15).Value = "-"
$wsData.Range($wsData.Cells.Item(710, 1), $wsData.Cells.Item(710, 15)).Style = "Normal"

# row 711
$wsData.Cells.Item(711, 1).Value = 100
$wsData.Cells.Item(711, 2).Value = 5
$wsData.Cells.Item(711, 3).Value = 0.5
$wsData.Cells.Item(711, 4).Value = "ba-cycle"
$wsData.Cells.Item(711, 5).Value = "2021-07-07 21:57:10.072471"
$wsData.Cells.Item(711, 6).Value = "-"
$wsData.Cells.Item(711, 7).Value = "-"
$wsData.Cells.Item(711, 8).Value = "-"
$wsData.Cells.Item(711, 9).Value = "-"
$wsData.Cells.Item(711, 10).Value = "-"
$wsData.Cells.Item(711, 11).Value = "-"
$wsData.Cells.Item(711, 12).Value = "-"
$wsData.Cells.Item(711, 13).Value = 5
$wsData.Cells.Item(711, 14).Value = 5
$wsData.Cells.Item(711, 15).Value = 5
$wsData.Range($wsData.Cells.Item(711, 1), $wsData.Cells.Item(711, 15)).Style = "Normal"

# row 712
$wsData.Cells.Item(712, 1).Value = 100
$wsData.Cells.Item(712, 2).Value = 5
$wsData.Cells.Item(712, 3).Value = 0.5
$wsData.Cells.Item(712, 4).Value = "er-no-cycle"
$wsData.Cells.Item(712, 5).Value = "2021-07-07 21:57:11.260660"
$wsData.Cells.Item(712, 6).Value = 13
$wsData.Cells.Item(712, 7).Value = 5
$wsData.Cells.Item(712, 8).Value = 5
$wsData.Cells.Item(712, 9).Value = 5
$wsData.Cells.Item(712, 10).Value = 7
$wsData.Cells.Item(712, 11).Value = 7
$wsData.Cells.Item(712, 12).Value = 7
$wsData.Cells.Item(712, 13).Value = "-"
$wsData.Cells.Item(712, 14).Value = "-"
$wsData.Cells.Item(712, 15).Value = "-"
$wsData.Range($wsData.Cells.Item(712, 1), $wsData.Cells.Item(712, 15)).Style = "Normal"

# row 713
$wsData.Cells.Item(713, 1).Value = 100
$wsData.Cells.Item(713, 2).Value = 5
$wsData.Cells.Item(713, 3).Value = 0.5
$wsData.Cells.Item(713, 4).Value = "er-cycle"
$wsData.Cells.Item(713, 5).Value = "2021-07-07 21:57:12.452814"
$wsData.Cells.Item(713, 6).Value = "-"
$wsData.Cells.Item(713, 7).Value = "-"
$wsData.Cells.Item(713, 8).Value = "-"
$wsData.Cells.Item(713, 9).Value = "-"
$wsData.Cells.Item(713, 10).Value = "-"
$wsData.Cells.Item(713, 11).Value = "-"
$wsData.Cells.Item(713, 12).Value = "-"
$wsData.Cells.Item(713, 13).Value = 7
$wsData.Cells.Item(713, 14).Value = 7
$wsData.Cells.Item(713, 15).Value = 7
$wsData.Range($wsData.Cells.Item(713, 1), $wsData.Cells.Item(713, 15)).Style = "Normal"

# row 714
$wsData.Cells.Item(714, 1).Value = 100
$wsData.Cells.Item(714, 2).Value = 5
$wsData.Cells.Item(714, 3).Value = 0.5
$wsData.Cells.Item(714, 4).Value = "ws-no-cycle"
$wsData.Cells.Item(714, 5).Value = "2021-07-07 21:57:13.604779"
$wsData.Cells.Item(714, 6).Value = 5
$wsData.Cells.Item(714, 7).Value = -3
$wsData.Cells.Item(714, 8).Value = 0
$wsData.Cells.Item(714, 9).Value = 1
$wsData.Cells.Item(714, 10).Value = 1
$wsData.Cells.Item(714, 11).Value = 1
$wsData.Cells.Item(714, 12).Value = 1
$wsData.Cells.Item(714, 13).Value = "-"
$wsData.Cells.Item(714, 14).Value = "-"
$wsData.Cells.Item(714, 15).Value = "-"
$wsData.Range($wsData.Cells.Item(714, 1), $wsData.Cells.Item(714, 15)).Style = "Normal"

# row 715
$wsData.Cells.Item(715, 1).Value = 100
$wsData.Cells.Item(715, 2).Value = 5
$wsData.Cells.Item(715, 3).Value = 0.5
$wsData.Cells.Item(715, 4).Value = "ws-cycle"
$wsData.Cells.Item(715, 5).Value = "2021-07-07 21:57:14.776810"
$wsData.Cells.Item(715, 6).Value = "-"
$wsData.Cells.Item(715, 7).Value = "-"
$wsData.Cells.Item(715, 8).Value = "-"
$wsData.Cells.Item(715, 9).Value = "-"
$wsData.Cells.Item(715, 10).Value = "-"
$wsData.Cells.Item(715, 11).Value = "-"
$wsData.Cells.Item(715, 12).Value = "-"
$wsData.Cells.Item(715, 13).Value = 1
$wsData.Cells.Item(715, 14).Value = 1
$wsData.Cells.Item(715, 15).Value = 1
$wsData.Range($wsData.Cells.Item(715, 1), $wsData.Cells.Item(715, 15)).Style = "Normal"

# row 716
$wsData.Cells.Item(716, 1).Value = 100
$wsData.Cells.Item(716, 2).Value = 5
$wsData.Cells.Item(716, 3).Value = 0.5
$wsData.Cells.Item(716, 4).Value = "cluster no cycle"
$wsData.Cells.Item(716, 5).Value = "2021-07-07 21:57:16.084809"
$wsData.Cells.Item(716, 6).Value = 5
$wsData.Cells.Item(716, 7).Value = 2
$wsData.Cells.Item(716, 8).Value = 4
$wsData.Cells.Item(716, 9).Value = 5
$wsData.Cells.Item(716, 10).Value = 5
$wsData.Cells.Item(716, 11).Value = 4
$wsData.Cells.Item(716, 12).Value = 4
$wsData.Cells.Item(716, 13).Value = "-"
$wsData.Cells.Item(716, 14).Value = "-"
$wsData.Cells.Item(716, 15).Value = "-"
$wsData.Range($wsData.Cells.Item(716, 1), $wsData.Cells.Item(716, 15)).Style = "Normal"

# row 717
$wsData.Cells.Item(717, 1).Value = 100
$wsData.Cells.Item(717, 2).Value = 5
$wsData.Cells.Item(717, 3).Value = 0.5
$wsData.Cells.Item(717, 4).Value = "cluster cycle"
$wsData.Cells.Item(717, 5).Value = "2021-07-07 21:57:17.268784"
$wsData.Cells.Item(717, 6).Value = "-"
$wsData.Cells.Item(717, 7).Value = "-"
$wsData.Cells.Item(717, 8).Value = "-"
$wsData.Cells.Item(717, 9).Value = "-"
$wsData.Cells.Item(717, 10).Value = "-"
$wsData.Cells.Item(717, 11).Value = "-"
$wsData.Cells.Item(717, 12).Value = "-"
$wsData.Cells.Item(717, 13).Value = 14
$wsData.Cells.Item(717, 14).Value = 14
$wsData.Cells.Item(717, 15).Value = 14
$wsData.Range($wsData.Cells.Item(717, 1), $wsData.Cells.Item(717, 15)).Style = "Normal"

# row 718
$wsData.Cells.Item(718, 1).Value = 100
$wsData.Cells.Item(718, 2).Value = 5
$wsData.Cells.Item(718, 3).Value = 0.5
$wsData.Cells.Item(718, 4).Value = "er"
$wsData.Cells.Item(718, 5).Value = "2021-07-07 22:25:53.530205"
$wsData.Cells.Item(718, 6).Value = 6
$wsData.Cells.Item(718, 7).Value = 3
$wsData.Cells.Item(718, 8).Value = 3
$wsData.Cells.Item(718, 9).Value = 3
$wsData.Cells.Item(718, 10).Value = 4
$wsData.Cells.Item(718, 11).Value = 0
$wsData.Cells.Item(718, 12).Value = 4
$wsData.Cells.Item(718, 13).Value = "-"
$wsData.Cells.Item(718, 14).Value = "-"
$wsData.Cells.Item(718, 15).Value = "-"
$wsData.Range($wsData.Cells.Item(718, 1), $wsData.Cells.Item(718, 15)).Style = "Normal"

# row 719
$wsData.Cells.Item(719, 1).Value = 100
$wsData.Cells.Item(719, 2).Value = 5
$wsData.Cells.Item(719, 3).Value = 0.5
$wsData.Cells.Item(719, 4).Value = "ws"
$wsData.Cells.Item(719, 5).Value = "2021-07-07 22:29:56.838401"
$wsData.Cells.Item(719, 6).Value = 18
$wsData.Cells.Item(719, 7).Value = 15
$wsData.Cells.Item(719, 8).Value = 15
$wsData.Cells.Item(719, 9).Value = 15
$wsData.Cells.Item(719, 10).Value = 17
$wsData.Cells.Item(719, 11).Value = 17
$wsData.Cells.Item(719, 12).Value = 17
$wsData.Cells.Item(719, 13).Value = "-"
$wsData.Cells.Item(719, 14).Value = "-"
$wsData.Cells.Item(719, 15).Value = "-"
$wsData.Range($wsData.Cells.Item(719, 1), $wsData.Cells.Item(719, 15)).Style = "Normal"

# row 720
$wsData.Cells.Item(720, 1).Value = 100
$wsData.Cells.Item(720, 2).Value = 5
$wsData.Cells.Item(720, 3).Value = 0.5
$wsData.Cells.Item(720, 4).Value = "ba"
$wsData.Cells.Item(720, 5).Value = "2021-07-07 22:31:40.806195"
$wsData.Cells.Item(720, 6).Value = 9
$wsData.Cells.Item(720, 7).Value = 6
$wsData.Cells.Item(720, 8).Value = 6
$wsData.Cells.Item(720, 9).Value = 6
$wsData.Cells.Item(720, 10).Value = 6
$wsData.Cells.Item(720, 11).Value = 6
$wsData.Cells.Item(720, 12).Value = 6
$wsData.Cells.Item(720, 13).Value = "-"
$wsData.Cells.Item(720, 14).Value = "-"
$wsData.Cells.Item(720, 15).Value = "-"
$wsData.Range($wsData.Cells.Item(720, 1), $wsData.Cells.Item(720, 15)).Style = "Normal"

# --- Runtimes sheet: new rows 711-721 ---
# row 711
$wsRuntimes.Cells.Item(711, 1).Value = 100
$wsRuntimes.Cells.Item(711, 2).Value = 5
$wsRuntimes.Cells.Item(711, 3).Value = 0.5
$wsRuntimes.Cells.Item(711, 4).Value = "ba-no-cycle"
$wsRuntimes.Cells.Item(711, 5).Value = "2021-07-07 21:57:08.700503"
$wsRuntimes.Cells.Item(711, 6).Value = 0.0007806000000001312
$wsRuntimes.Cells.Item(711, 7).Value = [double]"4.589999999993211e-05"
$wsRuntimes.Cells.Item(711, 8).Value = 0.01034859999999993
$wsRuntimes.Cells.Item(711, 9).Value = 0.06480240000000004
$wsRuntimes.Cells.Item(711, 10).Value = 0.03814990000000007
$wsRuntimes.Cells.Item(711, 11).Value = 0.002122400000000191
$wsRuntimes.Cells.Item(711, 12).Value = 0.00564349999999969
$wsRuntimes.Range($wsRuntimes.Cells.Item(711, 1), $wsRuntimes.Cells.Item(711, 12)).Style = "Normal"

# row 712
$wsRuntimes.Cells.Item(712, 1).Value = 100
$wsRuntimes.Cells.Item(712, 2).Value = 5
$wsRuntimes.Cells.Item(712, 3).Value = 0.5
$wsRuntimes.Cells.Item(712, 4).Value = "ba-cycle"
$wsRuntimes.Cells.Item(712, 5).Value = "2021-07-07 21:57:10.072471"
$wsRuntimes.Cells.Item(712, 6).Value = "-"
$wsRuntimes.Cells.Item(712, 7).Value = "-"
$wsRuntimes.Cells.Item(712, 8).Value = "-"
$wsRuntimes.Cells.Item(712, 9).Value = "-"
$wsRuntimes.Cells.Item(712, 10).Value = 0.0356285999999999
$wsRuntimes.Cells.Item(712, 11).Value = 0.002568799999999705
$wsRuntimes.Cells.Item(712, 12).Value = 0.005386599999999575
$wsRuntimes.Range($wsRuntimes.Cells.Item(712, 1), $wsRuntimes.Cells.Item(712, 12)).Style = "Normal"

# row 713
$wsRuntimes.Cells.Item(713, 1).Value = 100
$wsRuntimes.Cells.Item(713, 2).Value = 5
$wsRuntimes.Cells.Item(713, 3).Value = 0.5
$wsRuntimes.Cells.Item(713, 4).Value = "er-no-cycle"
$wsRuntimes.Cells.Item(713, 5).Value = "2021-07-07 21:57:11.260660"
$wsRuntimes.Cells.Item(713, 6).Value = 0.0004996000000003775
$wsRuntimes.Cells.Item(713, 7).Value = [double]"6.180000000011177e-05"
$wsRuntimes.Cells.Item(713, 8).Value = 0.0007615000000003036
$wsRuntimes.Cells.Item(713, 9).Value = 0.03180160000000054
$wsRuntimes.Cells.Item(713, 10).Value = 0.03520939999999939
$wsRuntimes.Cells.Item(713, 11).Value = 0.001506000000000007
$wsRuntimes.Cells.Item(713, 12).Value = 0.004215000000000302
$wsRuntimes.Range($wsRuntimes.Cells.Item(713, 1), $wsRuntimes.Cells.Item(713, 12)).Style = "Normal"

# row 714
$wsRuntimes.Cells.Item(714, 1).Value = 100
$wsRuntimes.Cells.Item(714, 2).Value = 5
$wsRuntimes.Cells.Item(714, 3).Value = 0.5
$wsRuntimes.Cells.Item(714, 4).Value = "er-cycle"
$wsRuntimes.Cells.Item(714, 5).Value = "2021-07-07 21:57:12.452814"
$wsRuntimes.Cells.Item(714, 6).Value = "-"
$wsRuntimes.Cells.Item(714, 7).Value = "-"
$wsRuntimes.Cells.Item(714, 8).Value = "-"
$wsRuntimes.Cells.Item(714, 9).Value = "-"
$wsRuntimes.Cells.Item(714, 10).Value = 0.03518810000000006
$wsRuntimes.Cells.Item(714, 11).Value = 0.001177300000000159
$wsRuntimes.Cells.Item(714, 12).Value = 0.003954600000000141
$wsRuntimes.Range($wsRuntimes.Cells.Item(714, 1), $wsRuntimes.Cells.Item(714, 12)).Style = "Normal"

# row 715
$wsRuntimes.Cells.Item(715, 1).Value = 100
$wsRuntimes.Cells.Item(715, 2).Value = 5
$wsRuntimes.Cells.Item(715, 3).Value = 0.5
$wsRuntimes.Cells.Item(715, 4).Value = "ws-no-cycle"
$wsRuntimes.Cells.Item(715, 5).Value = "2021-07-07 21:57:13.604779"
$wsRuntimes.Cells.Item(715, 6).Value = 0.0009300999999997117
$wsRuntimes.Cells.Item(715, 7).Value = [double]"4.490000000068051e-05"
$wsRuntimes.Cells.Item(715, 8).Value = 0.001740199999998637
$wsRuntimes.Cells.Item(715, 9).Value = 0.03593919999999962
$wsRuntimes.Cells.Item(715, 10).Value = 0.0347767000000001
$wsRuntimes.Cells.Item(715, 11).Value = 0.002935799999999489
$wsRuntimes.Cells.Item(715, 12).Value = 0.004860799999999443
$wsRuntimes.Range($wsRuntimes.Cells.Item(715, 1), $wsRuntimes.Cells.Item(715, 12)).Style = "Normal"

# row 716
$wsRuntimes.Cells.Item(716, 1).Value = 100
$wsRuntimes.Cells.Item(716, 2).Value = 5
$wsRuntimes.Cells.Item(716, 3).Value = 0.5
$wsRuntimes.Cells.Item(716, 4).Value = "ws-cycle"
$wsRuntimes.Cells.Item(716, 5).Value = "2021-07-07 21:57:14.776810"
$wsRuntimes.Cells.Item(716, 6).Value = "-"
$wsRuntimes.Cells.Item(716, 7).Value = "-"
$wsRuntimes.Cells.Item(716, 8).Value = "-"
$wsRuntimes.Cells.Item(716, 9).Value = "-"
$wsRuntimes.Cells.Item(716, 10).Value = 0.03404320000000105
$wsRuntimes.Cells.Item(716, 11).Value = 0.00248810000000077
$wsRuntimes.Cells.Item(716, 12).Value = 0.004376300000000555
$wsRuntimes.Range($wsRuntimes.Cells.Item(716, 1), $wsRuntimes.Cells.Item(716, 12)).Style = "Normal"

# row 717
$wsRuntimes.Cells.Item(717, 1).Value = 100
$wsRuntimes.Cells.Item(717, 2).Value = 5
$wsRuntimes.Cells.Item(717, 3).Value = 0.5
$wsRuntimes.Cells.Item(717, 4).Value = "cluster no cycle"
$wsRuntimes.Cells.Item(717, 5).Value = "2021-07-07 21:57:16.084809"
$wsRuntimes.Cells.Item(717, 6).Value = 0.001307899999998696
$wsRuntimes.Cells.Item(717, 7).Value = [double]"6.310000000020466e-05"
$wsRuntimes.Cells.Item(717, 8).Value = 0.1136543999999997
$wsRuntimes.Cells.Item(717, 9).Value = 0.03617579999999876
$wsRuntimes.Cells.Item(717, 10).Value = 0.03291789999999928
$wsRuntimes.Cells.Item(717, 11).Value = 0.001879999999999882
$wsRuntimes.Cells.Item(717, 12).Value = 0.007716600000000184
$wsRuntimes.Range($wsRuntimes.Cells.Item(717, 1), $wsRuntimes.Cells.Item(717, 12)).Style = "Normal"

# row 718
$wsRuntimes.Cells.Item(718, 1).Value = 100
$wsRuntimes.Cells.Item(718, 2).Value = 5
$wsRuntimes.Cells.Item(718, 3).Value = 0.5
$wsRuntimes.Cells.Item(718, 4).Value = "cluster cycle"
$wsRuntimes.Cells.Item(718, 5).Value = "2021-07-07 21:57:17.268784"
$wsRuntimes.Cells.Item(718, 6).Value = "-"
$wsRuntimes.Cells.Item(718, 7).Value = "-"
$wsRuntimes.Cells.Item(718, 8).Value = "-"
$wsRuntimes.Cells.Item(718, 9).Value = "-"
$wsRuntimes.Cells.Item(718, 10).Value = 0.03537269999999992
$wsRuntimes.Cells.Item(718, 11).Value = 0.003301100000001611
$wsRuntimes.Cells.Item(718, 12).Value = 0.008220299999999625
$wsRuntimes.Range($wsRuntimes.Cells.Item(718, 1), $wsRuntimes.Cells.Item(718, 12)).Style = "Normal"

# row 719
$wsRuntimes.Cells.Item(719, 1).Value = 100
$wsRuntimes.Cells.Item(719, 2).Value = 5
$wsRuntimes.Cells.Item(719, 3).Value = 0.5
$wsRuntimes.Cells.Item(719, 4).Value = "er"
$wsRuntimes.Cells.Item(719, 5).Value = "2021-07-07 22:25:53.530205"
$wsRuntimes.Cells.Item(719, 6).Value = 0.0006656999999998803
$wsRuntimes.Cells.Item(719, 7).Value = [double]"9.639999999988547e-05"
$wsRuntimes.Cells.Item(719, 8).Value = 0.001335700000000051
$wsRuntimes.Cells.Item(719, 9).Value = 0.05504299999999995
$wsRuntimes.Cells.Item(719, 10).Value = 0.05082720000000007
$wsRuntimes.Cells.Item(719, 11).Value = 0.0032582000000001
$wsRuntimes.Cells.Item(719, 12).Value = 0.008643699999999921
$wsRuntimes.Range($wsRuntimes.Cells.Item(719, 1), $wsRuntimes.Cells.Item(719, 12)).Style = "Normal"

# row 720
$wsRuntimes.Cells.Item(720, 1).Value = 100
$wsRuntimes.Cells.Item(720, 2).Value = 5
$wsRuntimes.Cells.Item(720, 3).Value = 0.5
$wsRuntimes.Cells.Item(720, 4).Value = "ws"
$wsRuntimes.Cells.Item(720, 5).Value = "2021-07-07 22:29:56.838401"
$wsRuntimes.Cells.Item(720, 6).Value = 0.0007159999999999389
$wsRuntimes.Cells.Item(720, 7).Value = [double]"6.680000000014452e-05"
$wsRuntimes.Cells.Item(720, 8).Value = 0.004442199999999952
$wsRuntimes.Cells.Item(720, 9).Value = 0.05318109999999998
$wsRuntimes.Cells.Item(720, 10).Value = 0.06027199999999988
$wsRuntimes.Cells.Item(720, 11).Value = 0.003525299999999953
$wsRuntimes.Cells.Item(720, 12).Value = 0.01144409999999985
$wsRuntimes.Range($wsRuntimes.Cells.Item(720, 1), $wsRuntimes.Cells.Item(720, 12)).Style = "Normal"

# row 721
$wsRuntimes.Cells.Item(721, 1).Value = 100
$wsRuntimes.Cells.Item(721, 2).Value = 5
$wsRuntimes.Cells.Item(721, 3).Value = 0.5
$wsRuntimes.Cells.Item(721, 4).Value = "ba"
$wsRuntimes.Cells.Item(721, 5).Value = "2021-07-07 22:31:40.806195"
$wsRuntimes.Cells.Item(721, 6).Value = 0.0007654000000001382
$wsRuntimes.Cells.Item(721, 7).Value = [double]"7.099999999993223e-05"
$wsRuntimes.Cells.Item(721, 8).Value = 0.003318500000000002
$wsRuntimes.Cells.Item(721, 9).Value = 0.0458464999999999
$wsRuntimes.Cells.Item(721, 10).Value = 0.05280600000000013
$wsRuntimes.Cells.Item(721, 11).Value = 0.002782000000000062
$wsRuntimes.Cells.Item(721, 12).Value = 0.009597299999999809
$wsRuntimes.Range($wsRuntimes.Cells.Item(721, 1), $wsRuntimes.Cells.Item(721, 12)).Style = "Normal"

Write-Host "Added 11 rows to Data sheet (710-720) and 11 rows to Runtimes sheet (711-721)"
